$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1094.2727
$ws.Range("I2").Value = 1190.7
$ws.Range("J2").Value = 130
$ws.Range("K2").Value = 1190.7
$ws.Range("L2").Value = 130
$ws.Range("M2").Value = -1077.7
$ws.Range("N2").Value = -356
$ws.Range("H33").Value = 615.4583
$ws.Range("I33").Value = 456.72223
$ws.Range("J33").Value = 1091.6666
$ws.Range("K33").Value = 456.72223
$ws.Range("L33").Value = 1091.6666
$ws.Range("M33").Value = -227.72223
$ws.Range("N33").Value = -1549.6666
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 5265618.5
$ws.Range("I132").Value = 6453494
$ws.Range("J132").Value = 5027.143
$ws.Range("K132").Value = 19360482
$ws.Range("L132").Value = 15081.429
$ws.Range("M132").Value = -19357952
$ws.Range("N132").Value = -20141.429
$ws.Range("H137").Value = 1783.303
$ws.Range("I137").Value = 1197.1111
$ws.Range("J137").Value = 2486.7334
$ws.Range("K137").Value = 3591.3333
$ws.Range("L137").Value = 7460.2002
$ws.Range("M137").Value = -1041.3333
$ws.Range("N137").Value = -12560.2002
$ws.Range("H138").Value = 3668.1228
$ws.Range("I138").Value = 1812.1951
$ws.Range("J138").Value = 8423.9375
$ws.Range("K138").Value = 5436.5853
$ws.Range("L138").Value = 25271.8125
$ws.Range("M138").Value = -296.5852999999997
$ws.Range("N138").Value = -35551.8125
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4995.409
$ws.Range("I32").Value = 3563
$ws.Range("J32").Value = 11441.25
$ws.Range("K32").Value = 3563
$ws.Range("L32").Value = 11441.25
$ws.Range("M32").Value = -3276
$ws.Range("N32").Value = -12015.25
$ws.Range("H61").Value = 2765.0488
$ws.Range("I61").Value = 963.7368
$ws.Range("K61").Value = 963.7368
$ws.Range("M61").Value = -751.7368
$ws.Range("H132").Value = 38466836
$ws.Range("I132").Value = 62505910
$ws.Range("J132").Value = 4319.8
$ws.Range("K132").Value = 187517730
$ws.Range("L132").Value = 12959.4
$ws.Range("M132").Value = -187515200
$ws.Range("N132").Value = -18019.4
$ws.Range("H136").Value = 2765.0488
$ws.Range("I136").Value = 963.7368
$ws.Range("K136").Value = 2891.2104
$ws.Range("M136").Value = -341.2103999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 739.71875
$ws.Range("I80").Value = 747.8125
$ws.Range("J80").Value = 731.625
$ws.Range("K80").Value = 747.8125
$ws.Range("L80").Value = 731.625
$ws.Range("M80").Value = 250.1875
$ws.Range("N80").Value = -2727.625
$ws.Range("H83").Value = 739.71875
$ws.Range("I83").Value = 747.8125
$ws.Range("J83").Value = 731.625
$ws.Range("K83").Value = 3739.0625
$ws.Range("L83").Value = 3658.125
$ws.Range("M83").Value = 1252.9375
$ws.Range("N83").Value = -13642.125
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2517.2666
$ws.Range("I31").Value = 1483.7435
$ws.Range("J31").Value = 4436.6665
$ws.Range("K31").Value = 1483.7435
$ws.Range("L31").Value = 4436.6665
$ws.Range("M31").Value = -1188.7435
$ws.Range("N31").Value = -5026.6665
$ws.Range("H34").Value = 2517.2666
$ws.Range("I34").Value = 1483.7435
$ws.Range("J34").Value = 4436.6665
$ws.Range("K34").Value = 1483.7435
$ws.Range("L34").Value = 4436.6665
$ws.Range("M34").Value = -1281.7435
$ws.Range("N34").Value = -4840.6665
$ws.Range("H70").Value = 32166.666
$ws.Range("J70").Value = 32166.666
$ws.Range("L70").Value = 32166.666
$ws.Range("N70").Value = -32796.666
$ws.Range("H73").Value = 32166.666
$ws.Range("J73").Value = 32166.666
$ws.Range("L73").Value = 32166.666
$ws.Range("N73").Value = -34350.666
$ws.Range("H86").Value = 3460.8484
$ws.Range("I86").Value = 2406.2222
$ws.Range("J86").Value = 4726.4
$ws.Range("K86").Value = 2406.2222
$ws.Range("L86").Value = 4726.4
$ws.Range("M86").Value = -1283.2222
$ws.Range("N86").Value = -6972.4
$ws.Range("H89").Value = 3460.8484
$ws.Range("I89").Value = 2406.2222
$ws.Range("J89").Value = 4726.4
$ws.Range("K89").Value = 12031.111
$ws.Range("L89").Value = 23632
$ws.Range("M89").Value = -6415.111000000001
$ws.Range("N89").Value = -34864
$ws.Range("H134").Value = 1728.5518
$ws.Range("I134").Value = 1345.8541
$ws.Range("J134").Value = 3565.5
$ws.Range("K134").Value = 4037.5623
$ws.Range("L134").Value = 10696.5
$ws.Range("M134").Value = -1502.5623
$ws.Range("N134").Value = -15766.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1363.0918
$ws.Range("I68").Value = 695.48486
$ws.Range("J68").Value = 1702.0308
$ws.Range("K68").Value = 2086.45458
$ws.Range("L68").Value = 5106.0924
$ws.Range("M68").Value = -1275.45458
$ws.Range("N68").Value = -6728.0924
$ws.Range("H71").Value = 1363.0918
$ws.Range("I71").Value = 695.48486
$ws.Range("J71").Value = 1702.0308
$ws.Range("K71").Value = 6259.363740000001
$ws.Range("L71").Value = 15318.2772
$ws.Range("M71").Value = -2203.363740000001
$ws.Range("N71").Value = -23430.2772
$ws.Range("H131").Value = 1910.909
$ws.Range("J131").Value = 1853.4615
$ws.Range("L131").Value = 5560.3845
$ws.Range("N131").Value = -15640.3845
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3937.6924
$ws.Range("I80").Value = 3630
$ws.Range("J80").Value = 4296.6665
$ws.Range("K80").Value = 3630
$ws.Range("L80").Value = 4296.6665
$ws.Range("M80").Value = -2632
$ws.Range("N80").Value = -6292.6665
$ws.Range("H83").Value = 3937.6924
$ws.Range("I83").Value = 3630
$ws.Range("J83").Value = 4296.6665
$ws.Range("K83").Value = 18150
$ws.Range("L83").Value = 21483.3325
$ws.Range("M83").Value = -13158
$ws.Range("N83").Value = -31467.3325
$ws.Range("H132").Value = 3659.3684
$ws.Range("I132").Value = 2502.3076
$ws.Range("J132").Value = 6166.3335
$ws.Range("K132").Value = 7506.9228
$ws.Range("L132").Value = 18499.0005
$ws.Range("M132").Value = -4976.9228
$ws.Range("N132").Value = -23559.0005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3341.04
$ws.Range("I136").Value = 2970.0625
$ws.Range("J136").Value = 4000.5557
$ws.Range("K136").Value = 8910.1875
$ws.Range("L136").Value = 12001.6671
$ws.Range("M136").Value = -6360.1875
$ws.Range("N136").Value = -17101.6671
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 27333.334
$ws.Range("J75").Value = 27333.334
$ws.Range("L75").Value = 27333.334
$ws.Range("N75").Value = -29205.334
$ws.Range("H78").Value = 27333.334
$ws.Range("J78").Value = 27333.334
$ws.Range("L78").Value = 82000.00199999999
$ws.Range("N78").Value = -91360.00199999999
$ws.Range("H108").Value = 33000
$ws.Range("J108").Value = 33000
$ws.Range("L108").Value = 33000
$ws.Range("N108").Value = -40680
$ws.Range("H123").Value = 43493.332
$ws.Range("J123").Value = 43493.332
$ws.Range("L123").Value = 43493.332
$ws.Range("N123").Value = -53293.332
$ws.Range("H141").Value = 28541.666
$ws.Range("J141").Value = 28541.666
$ws.Range("L141").Value = 28541.666
$ws.Range("N141").Value = -38901.666
